# Upload-refresh of the day/place mapping table:
#  - the date series in A4:A13 is recomputed as plain literal values
#    (the incremental "+1" formulas used for some of the rows are dropped)
#  - the whole series is shifted forward by one day starting at A5
#  - the window selection/scroll moves from B5:B13 down to A3:A13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-bearing cells with their new literal date serials.
# Setting .Value (rather than .Formula) clears any existing formula so the
# cell becomes a plain number, matching the target workbook.
$ws.Range("A4").Value  = 45736
$ws.Range("A5").Value  = 45737
$ws.Range("A6").Value  = 45738
$ws.Range("A7").Value  = 45739
$ws.Range("A8").Value  = 45740
$ws.Range("A9").Value  = 45741
$ws.Range("A10").Value = 45742
$ws.Range("A11").Value = 45743
$ws.Range("A12").Value = 45744
$ws.Range("A13").Value = 45745

# Scroll the window so row 5 is at the top of the viewport.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1

# Update the active selection to A3:A13 (active cell A3), replacing the
# previous B5:B13 selection.
$ws.Range("A3:A13").Select()
